# Apply updated correlation-analysis numbers ("made the main one be max")
# Sheet "all_tools" and sheet "openjml" get refreshed warning counts (F, G)
# and recomputed correlation statistics (I, J, K, L) for a handful of rows.

$wb = $excel.ActiveWorkbook

# ---- Sheet: all_tools ----
$ws = $wb.Worksheets.Item("all_tools")

# Row 9
$ws.Range("F9").Value = 94
$ws.Range("G9").Value = 1200
$ws.Range("I9").Value = -0.1679930849817192
$ws.Range("J9").Value = 0.01571209940531982
$ws.Range("K9").Value = -0.2581562594294852
$ws.Range("L9").Value = 0.009508181415596517

# Row 10
$ws.Range("G10").Value = 863
$ws.Range("I10").Value = -0.01788960397609135
$ws.Range("J10").Value = 0.8831700141519032
$ws.Range("K10").Value = -0.01605248793883186
$ws.Range("L10").Value = 0.9118983828429555

# Row 11
$ws.Range("G11").Value = 863
$ws.Range("I11").Value = -0.0319066731228802
$ws.Range("J11").Value = 0.7837672300250216
$ws.Range("K11").Value = -0.03045075791865187
$ws.Range("L11").Value = 0.8337285990222609

# Row 12
$ws.Range("G12").Value = 863
$ws.Range("I12").Value = -0.2441570787714312
$ws.Range("J12").Value = 0.03455416553919029
$ws.Range("K12").Value = -0.3094340106343606
$ws.Range("L12").Value = 0.02876701130235894

# ---- Sheet: openjml ----
$ws2 = $wb.Worksheets.Item("openjml")

# Row 9
$ws2.Range("F9").Value = 69
$ws2.Range("G9").Value = 808
$ws2.Range("I9").Value = -0.1508932577851104
$ws2.Range("J9").Value = 0.03976539977998658
$ws2.Range("K9").Value = -0.2101007063448729
$ws2.Range("L9").Value = 0.03589890255176344

# Row 10
$ws2.Range("F10").Value = 41
$ws2.Range("G10").Value = 219

# Row 11
$ws2.Range("F11").Value = 41
$ws2.Range("G11").Value = 219

# Row 12
$ws2.Range("F12").Value = 41
$ws2.Range("G12").Value = 219
